# Update the "final tweaks to slides" date placeholders:
#   5/26/2025                -> 6/5/2025                 (short numeric date fields)
#   Monday, May 26, 2025     -> Thursday, June 5, 2025    (long-form date fields)
#
# These cached date/time fields live on the slide master, every slide
# layout, the handout master and the notes master (not on the slides
# themselves).

$p = $ppt.ActivePresentation

$longOld  = "Monday, May 26, 2025"
$longNew  = "Thursday, June 5, 2025"
$shortOld = "5/26/2025"
$shortNew = "6/5/2025"

function Update-DateField($shape, [string]$newText) {
    $shape.TextFrame.TextRange.Text = $newText
}

# --- Slide Master -----------------------------------------------------
$master = $p.SlideMaster
Update-DateField $master.Shapes.Item(3) $longNew

# --- Slide Layouts ------------------------------------------------------
$layoutDateIdx = @{
    1 = 3
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 3
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $idx = $layoutDateIdx[$i]
    Update-DateField $layout.Shapes.Item($idx) $longNew
}

# --- Handout Master ------------------------------------------------------
$handout = $p.HandoutMaster
Update-DateField $handout.Shapes.Item(2) $shortNew

# --- Notes Master ------------------------------------------------------
$notes = $p.NotesMaster
Update-DateField $notes.Shapes.Item(2) $shortNew
